# Add new power plant types (rows 19-24) to the DPbES (Dispatch Priority by
# Electricity Source) subscript, each with a dispatch priority value of 3
# across all years (columns B:AK), mirroring issue #280 / #99.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")

$newSources = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$startRow = 19
for ($i = 0; $i -lt $newSources.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newSources[$i]
    for ($col = 2; $col -le 37; $col++) {
        $ws.Cells.Item($row, $col).Value = 3
    }
}

# Leave the final selection on the cell below the newly entered data, as the
# original author's session did, without disturbing which sheet tab is
# active in the workbook.
$origActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("A25").Select()
$origActive.Activate()
